# RPA datasets push 2023-11-08
# Insert a new top data row ("쏘닉스") above the existing listing table,
# shifting all prior data rows down by one, and update the subscription
# date (column A) of the row that used to be first (KB제27호스팩).
#
# The sheet is a plain data dump (no formulas), so the edit is a row
# insert + value fill. We avoid Range.Value/Value2 direct assignment of
# date-looking text (e.g. "2023-10-26") because Excel's type inference
# would silently convert it to a real date serial (and allocate a new
# number-format style) instead of leaving it as text like the source
# data. Routing everything through Copy + PasteSpecial(xlPasteValues)
# keeps cells text-typed and avoids touching cell styles.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$xlPasteValues = -4163

# 1) Shift existing data rows 2-17 down to 3-18.
$ws.Range("A2:T17").Copy()
$ws.Range("A3").PasteSpecial($xlPasteValues)
$excel.CutCopyMode = $false

# 2) Stage the brand-new row (쏘닉스) far below the table so typing
#    date-like text there doesn't touch the live row-2 styles; values
#    that look like dates are apostrophe-prefixed to force text entry.
$staging = 30

$ws.Cells.Item($staging, 1).Value2  = "'2023-10-26"      # A 청약일
$ws.Cells.Item($staging, 2).Value2  = "쏘닉스"             # B 회사명
$ws.Cells.Item($staging, 3).Value2  = "KB"                # C 대표주관회사
$ws.Cells.Item($staging, 4).Value2  = "'2023-10-31"       # D 납입일
$ws.Cells.Item($staging, 5).Value2  = "'2023-11-07"       # E 상장일
$ws.Cells.Item($staging, 6).Value2  = 27000000            # F 공모금액(천원)
$ws.Cells.Item($staging, 7).Value2  = 3600000             # G 공모주수
$ws.Cells.Item($staging, 8).Value2  = "-"                 # H 기준가(평가가치)
$ws.Cells.Item($staging, 9).Value2  = 5000                # I 1차발행가액(하단)
$ws.Cells.Item($staging, 10).Value2 = 7000                # J 1차발행가액(상단)
$ws.Cells.Item($staging, 11).Value2 = "-"                 # K 수요예측가중평균가
$ws.Cells.Item($staging, 12).Value2 = 7500                # L 확정발행가액
$ws.Cells.Item($staging, 13).Value2 = "-"                 # M 결정비율
$ws.Cells.Item($staging, 14).Value2 = "-"                 # N 공모비율
$ws.Cells.Item($staging, 15).Value2 = 0                   # O 구주매출비중
$ws.Cells.Item($staging, 16).Value2 = "-"                 # P 상장요건
$ws.Cells.Item($staging, 17).Value2 = "-"                 # Q 코넥스여부
$ws.Cells.Item($staging, 18).Value2 = "773.94 : 1"        # R 경쟁률
$ws.Cells.Item($staging, 19).Value2 = "-"                 # S 인수수수료(천원)
$ws.Cells.Item($staging, 20).Value2 = "-"                 # T 수수료율

$ws.Range($ws.Cells.Item($staging, 1), $ws.Cells.Item($staging, 20)).Copy()
$ws.Range("A2").PasteSpecial($xlPasteValues)
$excel.CutCopyMode = $false

$ws.Range($ws.Cells.Item($staging, 1), $ws.Cells.Item($staging, 20)).Clear()

# 3) The row that used to be first (KB제27호스팩, now row 3) had its
#    subscription date (청약일) refreshed from 2023-10-24 to 2023-10-26.
$ws.Cells.Item(3, 1).Value2 = "'2023-10-26"

# 4) Restore the rest of that row's original values (the copy-down in
#    step 1 already carried them, this just re-asserts column A was the
#    only change for that record, matching the source refresh).
